# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
$wb = $excel.ActiveWorkbook

# --- OFF sheet updates ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 454
$wsOff.Range("C2").Value = 326
$wsOff.Range("D2").Value = 91
$wsOff.Range("E2").Value = 41

# --- DEF sheet updates ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 481
$wsDef.Range("C2").Value = 321
$wsDef.Range("D2").Value = 116
$wsDef.Range("E2").Value = 57
$wsDef.Range("F2").Value = 10
